$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = -2
$ws.Range("F16").Value = 1
$ws.Range("F25").Value = 2
$ws.Range("F27").Value = 2
$ws.Range("F30").Value = -3
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = -3
$ws.Range("F33").Value = 3
$ws.Range("F39").Value = 3
$ws.Range("F40").Value = -1
$ws.Range("F46").Value = 4
$ws.Range("F50").Value = 4
$ws.Range("F54").Value = 1
$ws.Range("F60").Value = 3
$ws.Range("F61").Value = 0
$ws.Range("F64").Value = -3
$ws.Range("F69").Value = 3
